# Update odds for row 3 (MJCkWIE8 - Liverpool M. vs Wanderers) on the active sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value  = 2.1
$ws.Range("H3").Value  = 3.2
$ws.Range("I3").Value  = 3.2
$ws.Range("J3").Value  = 2.88
$ws.Range("M3").Value  = 1.07
$ws.Range("N3").Value  = 9
$ws.Range("O3").Value  = 1.33
$ws.Range("P3").Value  = 3.25
$ws.Range("Q3").Value  = 2.08
$ws.Range("R3").Value  = 1.73
$ws.Range("S3").Value  = 1.44
$ws.Range("T3").Value  = 2.63
$ws.Range("U3").Value  = 1.83
$ws.Range("V3").Value  = 1.83
$ws.Range("W3").Value  = 7
$ws.Range("Y3").Value  = 9.5
$ws.Range("Z3").Value  = 21
$ws.Range("AA3").Value = 19
$ws.Range("AC3").Value = 9
$ws.Range("AG3").Value = 9.5
$ws.Range("AJ3").Value = 34
$ws.Range("AL3").Value = 41
$ws.Range("AM3").Value = 301
$ws.Range("AN3").Value = 4.33
$ws.Range("AP3").Value = 23
$ws.Range("AR3").Value = 67
$ws.Range("AT3").Value = 2.63
$ws.Range("AW3").Value = 5
